$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45917
$ws.Range("B2").Value = 107.12
$ws.Range("C2").Value = 104.78
$ws.Range("D2").Value = 102.9
$ws.Range("E2").Value = 99.91
$ws.Range("F2").Value = 98.5
$ws.Range("G2").Value = 102
$ws.Range("H2").Value = 104.78
$ws.Range("I2").Value = 118.56
$ws.Range("J2").Value = 112.08
$ws.Range("K2").Value = 95.05
$ws.Range("L2").Value = 57.9
$ws.Range("M2").Value = 37.6
$ws.Range("N2").Value = 19.68
$ws.Range("O2").Value = 19.57
$ws.Range("P2").Value = 15.1
$ws.Range("Q2").Value = 16.4
$ws.Range("R2").Value = 20.33
$ws.Range("S2").Value = 50
$ws.Range("T2").Value = 83.12
$ws.Range("U2").Value = 111.76
$ws.Range("V2").Value = 150
$ws.Range("W2").Value = 215
$ws.Range("X2").Value = 117.11
$ws.Range("Y2").Value = 104.99
$ws.Range("Z2").Value = 86.01000000000001
$ws.Range("AB2").Value = 146.78
$ws.Range("AD2").Value = 182.5
$ws.Range("AE2").Value = "6h-8h"
$ws.Range("AF2").Value = 111.67
$ws.Range("AG2").Value = "10h-18h"
